# Continue building out the "find images" helper table on the stimuli sheet:
# fill in the "carrier" (column D) values for the practice/generic word rows,
# add the pair_kind (column J) markers for the video/audio pair rows, and
# populate the newly expanded unique_video / unique_audio block (rows 14-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5): carrier word used for each practice item.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic pair rows (6-9): mark which ones are the unique video/audio pairs.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: new unique_video / unique_audio kind + carrier entries.
$uniqueRows = @(
    @{ Row = 14; Kind = "unique_video"; Carrier = "look" },
    @{ Row = 15; Kind = "unique_video"; Carrier = "look" },
    @{ Row = 16; Kind = "unique_video"; Carrier = "where" },
    @{ Row = 17; Kind = "unique_video"; Carrier = "where" },
    @{ Row = 18; Kind = "unique_audio"; Carrier = "can" },
    @{ Row = 19; Kind = "unique_audio"; Carrier = "can" },
    @{ Row = 20; Kind = "unique_audio"; Carrier = "do" },
    @{ Row = 21; Kind = "unique_audio"; Carrier = "do" }
)

foreach ($item in $uniqueRows) {
    $ws.Cells.Item($item.Row, 3).Value = $item.Kind      # column C = kind
    $ws.Cells.Item($item.Row, 4).Value = $item.Carrier   # column D = carrier
}
